# Auto commit at 2025-09-07 7:37:50.25
# Update Metrics values (B2:B13), move the active sheet / selection from
# "today" to "Metrics", and bump the TODAY()-1 formula on "today" forward
# one day (its cached value recalculates automatically from the new
# system date, but we also let the engine recompute dependents of Metrics).

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# New values for Metrics!B2:B13
$metrics.Range("B2").Value  = 96622.579999999987
$metrics.Range("B3").Value  = 78580.48000000001
$metrics.Range("B4").Value  = 30709.390000000003
$metrics.Range("B5").Value  = 3790
$metrics.Range("B6").Value  = 4015873.4599999995
$metrics.Range("B7").Value  = 3406107.959999999
$metrics.Range("B8").Value  = 1160075.07
$metrics.Range("B9").Value  = 154950
$metrics.Range("B10").Value = 32481197.260999829
$metrics.Range("B11").Value = 19435978.030000005
$metrics.Range("B12").Value = 11441783.960000001
$metrics.Range("B13").Value = 1252577

# Recalculate so formulas on "today" (and elsewhere) pick up the new values.
$excel.Calculate()

# Move the selection / active-sheet state: "today" was the active tab with
# selection D5; now "Metrics" is active with selection F23, and "today"'s
# selection becomes I11.
$today.Range("I11").Select()
$metrics.Select()
$metrics.Range("F23").Select()

$wb.Save()
